$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header labels (row 1) for columns AB..AK -- shared-string additions
# ---------------------------------------------------------------------------
$headers = @("trf_ppn", "trf_pph", "trf_ppnbm", "trf_ppnbm_t", "trf_bmad", "trf_bmad_t", "trf_bk", "trf_bk_t", "bk_nilai_awal", "bk_nilai_akhir")
$cols = @("AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# Data rows 2 and 3, columns AB..AK
# ---------------------------------------------------------------------------
$values = @(10, 2.5, 0, 0, 0, 0, 0, 0, 0, 0)

for ($r = 2; $r -le 3; $r++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Column widths for the new columns
# ---------------------------------------------------------------------------
$ws.Columns("AD").ColumnWidth = 9.333333333333334
$ws.Columns("AE").ColumnWidth = 11.333333333333334
$ws.Columns("AG").ColumnWidth = 10.166666666666666
$ws.Columns("AH").ColumnWidth = 10.666666666666666
$ws.Columns("AI").ColumnWidth = 8.666666666666666
$ws.Columns("AJ").ColumnWidth = 14.333333333333334
$ws.Columns("AK").ColumnWidth = 12.666666666666666

# ---------------------------------------------------------------------------
# View: selection / top-left cell
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 1
$ws.Range("AJ7").Select()
